# Atualização de bases das ligas, do dia: 02-03-2024 às 08:34
#
# 1) Four pairs of existing match rows had their data reordered (the two
#    rows in each pair swap all of their match data while keeping their
#    own row position / running "id" index in column A untouched).
# 2) Five brand-new match rows are appended at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows {
    param([int]$Row1, [int]$Row2)

    # Columns B..AC (2..29) hold the match data; column A (the running
    # index) must stay put on each row.
    for ($col = 2; $col -le 29; $col++) {
        $cell1 = $ws.Cells.Item($Row1, $col)
        $cell2 = $ws.Cells.Item($Row2, $col)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}

Swap-MatchRows 104 105
Swap-MatchRows 107 108
Swap-MatchRows 143 145
Swap-MatchRows 148 149

# ---------------------------------------------------------------------
# Append new rows 176-180, matching formatting of the last existing row.
# ---------------------------------------------------------------------

function New-MatchRow {
    param([int]$Row)

    # Column A keeps the bold/border/centered style used by every id cell.
    $ws.Cells.Item(175, 1).Copy()
    $ws.Cells.Item($Row, 1).PasteSpecial(-4122)

    # Column E keeps the custom date/time number format.
    $ws.Cells.Item(175, 5).Copy()
    $ws.Cells.Item($Row, 5).PasteSpecial(-4122)
}

New-MatchRow 176
$ws.Cells.Item(176, 1).Value = 174
$ws.Cells.Item(176, 2).Value = 7897430
$ws.Cells.Item(176, 3).Value = "Bolivia Primera División"
$ws.Cells.Item(176, 4).Value = "Bolivia Apertura"
$ws.Cells.Item(176, 5).Value = 45352.875
$ws.Cells.Item(176, 6).Value = "Royal Pari FC"
$ws.Cells.Item(176, 7).Value = "Blooming"
$ws.Cells.Item(176, 8).Value = 1
$ws.Cells.Item(176, 9).Value = 2
$ws.Cells.Item(176, 10).Value = "A"
$ws.Cells.Item(176, 11).Value = 2.5
$ws.Cells.Item(176, 12).Value = 3.5
$ws.Cells.Item(176, 13).Value = 2.625
$ws.Cells.Item(176, 14).Value = 2.55
$ws.Cells.Item(176, 15).Value = 3.5
$ws.Cells.Item(176, 16).Value = 2.7
$ws.Cells.Item(176, 17).Value = 0
$ws.Cells.Item(176, 18).Value = 1.825
$ws.Cells.Item(176, 19).Value = 1.975
$ws.Cells.Item(176, 20).Value = 2.75
$ws.Cells.Item(176, 21).Value = 1.975
$ws.Cells.Item(176, 22).Value = 1.825
$ws.Cells.Item(176, 23).Value = -1
$ws.Cells.Item(176, 24).Value = -1
$ws.Cells.Item(176, 25).Value = 1.7
$ws.Cells.Item(176, 26).Value = -1
$ws.Cells.Item(176, 27).Value = 0.9750000000000001
$ws.Cells.Item(176, 28).Value = 0.4875
$ws.Cells.Item(176, 29).Value = -0.5

New-MatchRow 177
$ws.Cells.Item(177, 1).Value = 175
$ws.Cells.Item(177, 2).Value = 7897428
$ws.Cells.Item(177, 3).Value = "Bolivia Primera División"
$ws.Cells.Item(177, 4).Value = "Bolivia Apertura"
$ws.Cells.Item(177, 5).Value = 45353.66666666666
$ws.Cells.Item(177, 6).Value = "The Strongest"
$ws.Cells.Item(177, 7).Value = "Real Tomayapo"
$ws.Cells.Item(177, 11).Value = 1.363
$ws.Cells.Item(177, 12).Value = 4.5
$ws.Cells.Item(177, 13).Value = 7
$ws.Cells.Item(177, 14).Value = 1.083
$ws.Cells.Item(177, 15).Value = 11
$ws.Cells.Item(177, 16).Value = 23
$ws.Cells.Item(177, 17).Value = -2.75
$ws.Cells.Item(177, 18).Value = 1.925
$ws.Cells.Item(177, 19).Value = 1.875
$ws.Cells.Item(177, 20).Value = 3.75
$ws.Cells.Item(177, 21).Value = 1.9
$ws.Cells.Item(177, 22).Value = 1.9
$ws.Cells.Item(177, 23).Value = 0
$ws.Cells.Item(177, 24).Value = 0
$ws.Cells.Item(177, 25).Value = 0
$ws.Cells.Item(177, 26).Value = 0
$ws.Cells.Item(177, 27).Value = 0

New-MatchRow 178
$ws.Cells.Item(178, 1).Value = 176
$ws.Cells.Item(178, 2).Value = 7897433
$ws.Cells.Item(178, 3).Value = "Bolivia Primera División"
$ws.Cells.Item(178, 4).Value = "Bolivia Apertura"
$ws.Cells.Item(178, 5).Value = 45353.77083333334
$ws.Cells.Item(178, 6).Value = "Jorge Wilstermann"
$ws.Cells.Item(178, 7).Value = "Bolivar"
$ws.Cells.Item(178, 11).Value = 2.4
$ws.Cells.Item(178, 12).Value = 3.3
$ws.Cells.Item(178, 13).Value = 2.6
$ws.Cells.Item(178, 14).Value = 3.2
$ws.Cells.Item(178, 15).Value = 3.4
$ws.Cells.Item(178, 16).Value = 2.25
$ws.Cells.Item(178, 17).Value = 0.25
$ws.Cells.Item(178, 18).Value = 1.825
$ws.Cells.Item(178, 19).Value = 1.975
$ws.Cells.Item(178, 20).Value = 2.75
$ws.Cells.Item(178, 21).Value = 1.975
$ws.Cells.Item(178, 22).Value = 1.825
$ws.Cells.Item(178, 23).Value = 0
$ws.Cells.Item(178, 24).Value = 0
$ws.Cells.Item(178, 25).Value = 0
$ws.Cells.Item(178, 26).Value = 0
$ws.Cells.Item(178, 27).Value = 0

New-MatchRow 179
$ws.Cells.Item(179, 1).Value = 177
$ws.Cells.Item(179, 2).Value = 7897432
$ws.Cells.Item(179, 3).Value = "Bolivia Primera División"
$ws.Cells.Item(179, 4).Value = "Bolivia Apertura"
$ws.Cells.Item(179, 5).Value = 45353.875
$ws.Cells.Item(179, 6).Value = "Guabira"
$ws.Cells.Item(179, 7).Value = "Always Ready"
$ws.Cells.Item(179, 11).Value = 2.2
$ws.Cells.Item(179, 12).Value = 3.5
$ws.Cells.Item(179, 13).Value = 2.75
$ws.Cells.Item(179, 14).Value = 1.727
$ws.Cells.Item(179, 15).Value = 4
$ws.Cells.Item(179, 16).Value = 4.5
$ws.Cells.Item(179, 17).Value = -0.75
$ws.Cells.Item(179, 18).Value = 1.875
$ws.Cells.Item(179, 19).Value = 1.925
$ws.Cells.Item(179, 20).Value = 2.75
$ws.Cells.Item(179, 21).Value = 1.925
$ws.Cells.Item(179, 22).Value = 1.875
$ws.Cells.Item(179, 23).Value = 0
$ws.Cells.Item(179, 24).Value = 0
$ws.Cells.Item(179, 25).Value = 0
$ws.Cells.Item(179, 26).Value = 0
$ws.Cells.Item(179, 27).Value = 0

New-MatchRow 180
$ws.Cells.Item(180, 1).Value = 178
$ws.Cells.Item(180, 2).Value = 7897431
$ws.Cells.Item(180, 3).Value = "Bolivia Primera División"
$ws.Cells.Item(180, 4).Value = "Bolivia Apertura"
$ws.Cells.Item(180, 5).Value = 45354.66666666666
$ws.Cells.Item(180, 6).Value = "Nacional Potosi"
$ws.Cells.Item(180, 7).Value = "Club Aurora"
$ws.Cells.Item(180, 11).Value = 1.571
$ws.Cells.Item(180, 12).Value = 3.6
$ws.Cells.Item(180, 13).Value = 5.25
$ws.Cells.Item(180, 14).Value = 1.7
$ws.Cells.Item(180, 15).Value = 3.5
$ws.Cells.Item(180, 16).Value = 4.5
$ws.Cells.Item(180, 17).Value = -0.75
$ws.Cells.Item(180, 18).Value = 1.925
$ws.Cells.Item(180, 19).Value = 1.875
$ws.Cells.Item(180, 20).Value = 3
$ws.Cells.Item(180, 21).Value = 1.825
$ws.Cells.Item(180, 22).Value = 1.975
$ws.Cells.Item(180, 23).Value = 0
$ws.Cells.Item(180, 24).Value = 0
$ws.Cells.Item(180, 25).Value = 0
$ws.Cells.Item(180, 26).Value = 0
$ws.Cells.Item(180, 27).Value = 0
